# RDCC-5182 SRD file version validation changes
# Adds a new "VERSION" worksheet at the end of the workbook containing the
# file-version marker cells, and makes it the active sheet (which in turn
# clears the tabSelected flag previously set on Sheet1).

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet so it lands at the end
# of the tab order.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$versionSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$versionSheet.Name = "VERSION"

# Populate the version marker cells.
$versionSheet.Range("A6").Value = "File version"
$versionSheet.Range("B6").Value = "vx.xx"

# Mirror the authored selection/active-cell state on the new sheet.
$versionSheet.Range("B6").Select()
